# BusRaider V1.7 BOM - "Fixes and tidy up" edit
# U6 (74HC138 DIL decoder) was actually fitted/sourced as the HCT variant
# (74HCT138), and its supplier info is switched from LCSC to Farnell.
# U13/U10 (74HCT74) and U5 (74LVC07AD,118) supplier info also updated to
# Farnell (adding a Farnell part-number hyperlink for the 74HCT138 and the
# 74LVC07AD part).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Row 25 - U6: 74HC138 DIL -> 74HCT138 DIL (and part 74HC138 -> 74HCT138),
# supplier LCSC/C86611 -> Farnell/1470803
$ws.Range("B25").Value = "74HCT138 DIL"
$ws.Range("F25").Value = "74HCT138"
$ws.Range("H25").Value = "Farnell"
$ws.Hyperlinks.Add($ws.Range("I25"), "https://uk.farnell.com/texas-instruments/sn74hct138n/ic-decoder-demux/dp/1470803") | Out-Null
$ws.Range("I25").NumberFormat = "General"
$ws.Range("I25").Value = 1470803

# Row 26 - U13,U10: 74HCT74, supplier now explicitly Farnell
$ws.Range("H26").Value = "Farnell"

# Row 27 - U5: 74LVC07AD,118, supplier LCSC/C6049 -> Farnell/2463753
$ws.Range("H27").Value = "Farnell"
$ws.Hyperlinks.Add($ws.Range("I27"), "https://uk.farnell.com/nexperia/74lvc07ad-118/buffer-hex-non-inverting-tssop/dp/2463753?st=74lvc07") | Out-Null
$ws.Range("I27").NumberFormat = "General"
$ws.Range("I27").Value = 2463753

# Update the last active selection to match (cursor ended up on I27)
$ws.Range("I27").Select()

$wb.Save()
